# Update of ch2//chip-seq draft
# Adds a new "Sheet2" after the existing "Sheet1" with a small timepoint
# (t1/t2/t3) breakdown table + ratio formulas, and makes it the active sheet.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 (becomes active automatically).
$ws = $wb.Worksheets.Add($null, $sheet1)

# Row 1
$ws.Range("B1").Value = "t1"
$ws.Range("E1").Value = 782
$ws.Range("G1").Formula = "=E1/SUM(E1:E3)"

# Row 2
$ws.Range("B2").Value = "t1"
$ws.Range("C2").Value = "t2"
$ws.Range("E2").Value = 64

# Row 3
$ws.Range("B3").Value = "t1"
$ws.Range("C3").Value = "t2"
$ws.Range("D3").Value = "t3"
$ws.Range("E3").Value = 717

# Row 4
$ws.Range("C4").Value = "t2"
$ws.Range("E4").Value = 3056
$ws.Range("G4").Formula = "=E4/SUM(E2:E5)"

# Row 5
$ws.Range("C5").Value = "t2"
$ws.Range("D5").Value = "t3"
$ws.Range("E5").Value = 1409

# Row 6
$ws.Range("D6").Value = "t3"
$ws.Range("E6").Value = 2311
$ws.Range("G6").Formula = "=E6/SUM(E6,E5,E3)"

# Match the saved selection/active cell on the new sheet.
$ws.Range("G6").Select()
